$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) SEC_Comm: populate row 11 with the new "WIND_ON" commodity
# ---------------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item("SEC_Comm")
$wsComm.Range("B11").Value2 = "NRG"
$wsComm.Range("C11").Value2 = "WIND_ON"
$wsComm.Range("D11").Value2 = "Wind Onshore"
$wsComm.Range("E11").Value2 = "PJ"
$wsComm.Range("G11").Value2 = "SEASON"

# Match the formatting already used on the row above (row 10, same table)
$cols = @("B","C","D","E","F","G","H","I")
foreach ($col in $cols) {
    $src = $wsComm.Range($col + "10")
    $dst = $wsComm.Range($col + "11")
    $dst.NumberFormat = $src.NumberFormat
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.Interior.Color = $src.Interior.Color
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
}

# ---------------------------------------------------------------------------
# 2) SEC_Processes: insert a new technology row (MIN_EX_WIND_ON) right after
#    the existing MIN_EX_WIND row, keeping the closing border row below it.
# ---------------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("SEC_Processes")
$wsProc.Rows("13:14").Insert()

$wsProc.Range("B13").Value2 = "MIN"
$wsProc.Range("D13").Value2 = "MIN_EX_WIND_ON"
$wsProc.Range("E13").Value2 = "Wind mine"
$wsProc.Range("F13").Value2 = "PJ"
$wsProc.Range("G13").Value2 = "Pja"
$wsProc.Range("H13").Value2 = "SEASON"

# Row 13 gets the "inner" row styling (same as row 11, the row above the
# pair that used to close the table), row 14 becomes the new closing /
# thick-bottom-border row (same styling the old row 12 used to carry).
$cols = @("B","C","D","E","F","G","H","I","J")
foreach ($col in $cols) {
    $srcInner = $wsProc.Range($col + "11")
    $dstInner = $wsProc.Range($col + "13")
    $dstInner.NumberFormat = $srcInner.NumberFormat
    $dstInner.Font.Name = $srcInner.Font.Name
    $dstInner.Font.Size = $srcInner.Font.Size
    $dstInner.Font.Bold = $srcInner.Font.Bold
    $dstInner.Interior.Color = $srcInner.Interior.Color
    $dstInner.HorizontalAlignment = $srcInner.HorizontalAlignment
    $dstInner.VerticalAlignment = $srcInner.VerticalAlignment

    $srcBottom = $wsProc.Range($col + "12")
    $dstBottom = $wsProc.Range($col + "14")
    $dstBottom.NumberFormat = $srcBottom.NumberFormat
    $dstBottom.Font.Name = $srcBottom.Font.Name
    $dstBottom.Font.Size = $srcBottom.Font.Size
    $dstBottom.Font.Bold = $srcBottom.Font.Bold
    $dstBottom.Interior.Color = $srcBottom.Interior.Color
    $dstBottom.HorizontalAlignment = $srcBottom.HorizontalAlignment
    $dstBottom.VerticalAlignment = $srcBottom.VerticalAlignment
    $dstBottom.Borders.Item(9).LineStyle = $srcBottom.Borders.Item(9).LineStyle
    $dstBottom.Borders.Item(9).Weight = $srcBottom.Borders.Item(9).Weight
}

# Row 12 itself now sits "inside" the table (no longer the last row), so it
# loses its thick bottom border and picks up the plain inner-row style.
foreach ($col in $cols) {
    $srcInner = $wsProc.Range($col + "11")
    $dst12 = $wsProc.Range($col + "12")
    $dst12.NumberFormat = $srcInner.NumberFormat
    $dst12.Font.Name = $srcInner.Font.Name
    $dst12.Font.Size = $srcInner.Font.Size
    $dst12.Font.Bold = $srcInner.Font.Bold
    $dst12.Interior.Color = $srcInner.Interior.Color
    $dst12.HorizontalAlignment = $srcInner.HorizontalAlignment
    $dst12.VerticalAlignment = $srcInner.VerticalAlignment
}

# ---------------------------------------------------------------------------
# 3) MIN_IMP: tiny import bound for the existing WIND mine, plus a new row
#    describing the WIND_ON mine that references the rows created above.
# ---------------------------------------------------------------------------
$wsMin = $wb.Worksheets.Item("MIN_IMP")
$wsMin.Range("E10").Value2 = 0.001
$wsMin.Range("E10").NumberFormat = "0.000"

$wsMin.Range("B11").Formula = "=SEC_Processes!D13"
$wsMin.Range("C11").Formula = "=SEC_Comm!D11"
$wsMin.Range("D11").Formula = "=SEC_Comm!C11"
$wsMin.Range("E11").Value2 = 0.001
$wsMin.Range("E11").NumberFormat = "0.000"

$cols = @("B","C","D")
foreach ($col in $cols) {
    $src = $wsMin.Range($col + "10")
    $dst = $wsMin.Range($col + "11")
    $dst.NumberFormat = $src.NumberFormat
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.Interior.Color = $src.Interior.Color
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
}

# ---------------------------------------------------------------------------
# 4) View state: the user ended up on MIN_IMP, having last touched a handful
#    of cells on each sheet.
# ---------------------------------------------------------------------------
$wsComm.Range("E22").Select()
$wsProc.Range("I18").Select()
$wsMin.Range("F11").Select()

$wsPP = $wb.Worksheets.Item("PP")
$wsPP.Range("I16").Select()

$wsDMD = $wb.Worksheets.Item("DMD")
$wsDMD.Range("I41").Select()

$wsMin.Activate()
